# UPDATE: Item Icons added and filled
#
# Source data rows (String_Data sheet) for the "Lottery supply" and
# "Traffic supply" item ladders got reshuffled: a new "Crystal Ball" /
# "​水晶球​" item was inserted into the Lottery ladder (replacing
# "Rabbit's Foot" / "​兔脚", whose slot is now filled by "Dice"), and the
# Traffic ladder's LV2/LV3 rungs ("Jump Pad" / "Grappling Hook") were
# swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("String_Data")

# --- Traffic supply ladder: swap LV2 (row 38) and LV3 (row 39) ---
# Row 38 (TRAFFIC_SUPPLY_LV2): Jump Pad / <200b>跳板  ->  Grappling Hook / 钩爪
$ws.Range("B38").Value = "Grappling Hook"
$ws.Range("C38").Value = "钩爪"

# Row 39 (TRAFFIC_SUPPLY_LV3): Grappling Hook / 钩爪  ->  Jump Pad / <200b>跳板
$ws.Range("B39").Value = "Jump Pad"
$ws.Range("C39").Value = ([char]0x200B + "跳板")

# --- Lottery supply ladder: Rabbit's Foot removed, Dice shifts up, new Crystal Ball added ---
# Row 48 (LOTTERY_SUPPLY_LV2): Rabbit's Foot / <200b>兔脚  ->  Dice / <200b>骰子<200b>
$ws.Range("B48").Value = "Dice"
$ws.Range("C48").Value = ([char]0x200B + "骰子" + [char]0x200B)

# Row 49 (LOTTERY_SUPPLY_LV3): Dice / <200b>骰子<200b>  ->  Crystal Ball / <200b>水晶球<200b>
$ws.Range("B49").Value = "Crystal Ball"
$ws.Range("C49").Value = ([char]0x200B + "水晶球" + [char]0x200B)

# --- Restore the selection to match the author's final cursor position ---
$ws.Range("C51").Select()
